$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Ativacao date: 01/01/2012 -> 01/01/2021 ---
# Direct Value assignment would be auto-parsed as a date serial by Excel
# (it's an unambiguous valid date), which would change the cell's stored
# type/style. Instead, build the literal text via a helper formula cell and
# paste-special VALUES ONLY so the destination keeps its existing
# text-cell typing and style.
$ws.Range("Z1").Formula = "=""01/01/2021"""
$ws.Range("Z1").Copy()
$ws.Range("B8").PasteSpecial(-4163)
$ws.Range("C8").PasteSpecial(-4163)
$ws.Range("Z1").Clear()

# --- Docentes responsaveis ---
$ws.Range("B13").Value = "11079086 - Herlandí de Souza Andrade"
$ws.Range("C13").Value = "11079086 - Herlandí de Souza Andrade"

# --- Criterio ---
$ws.Range("B20").Value = "Média Aritmética das atividades avaliativas realizadas."
$ws.Range("C20").Value = "Média Aritmética das atividades avaliativas realizadas."

# --- Norma de recuperacao (trailing period removed) ---
$ws.Range("B21").Value = "Média aritmética da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recuperação"
$ws.Range("C21").Value = "Média aritmética da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recuperação"

# --- New requirement row (row 25), mirroring row 24's layout/style ---
$ws.Range("B24:C24").Copy()
$ws.Range("B25:C25").PasteSpecial(-4122)
$newReq = "LOQ4240 -  Administração e Organização II  (Requisito fraco)`n"
$ws.Cells.Item(25, 2).Formula = $newReq
$ws.Cells.Item(25, 3).Formula = $newReq
$ws.Rows.Item(25).RowHeight = 30
